$d = $word.ActiveDocument

# Locate the "Fix Procedures" table: it's the one whose second row still
# holds the leftover "x" / "y" scaffolding placeholders.
$t = $null
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $candidate = $d.Tables.Item($i)
    if ($candidate.Rows.Count -ge 2 -and $candidate.Rows.Item(2).Cells.Item(1).Range.Text.Substring(0,1) -eq "x" -and $candidate.Rows.Item(2).Cells.Item(1).Range.Text.Length -eq 3) {
        $t = $candidate
        break
    }
}
if ($t -eq $null) {
    $t = $d.Tables.Item(4)
}

# Helper: replace a paragraph's content with the given inner OOXML by
# wrapping it in the "WordOpenXML" pkg:package envelope and inserting it
# over the paragraph's own range (this fully replaces the paragraph,
# letting us control runs / lastRenderedPageBreak placement precisely).
function Set-ParaXml($para, $innerXml) {
    $rng = $para.Range
    $full = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body>" + $innerXml + "</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
    [void]$rng.InsertXML($full)
}

# 1. Remove the leftover scaffolding row ("x" / "y").
$t.Rows.Item(2).Delete()

# 2. The "Delete Empty File~" row (now row 2): drop its leading
#    lastRenderedPageBreak, and move the page-break marker into the
#    middle of the description text (splitting it into two runs).
$row2 = $t.Rows.Item(2)
$p1 = $row2.Cells.Item(1).Range.Paragraphs.Item(1)
Set-ParaXml $p1 "<w:p><w:pPr><w:pStyle w:val='Compact'/></w:pPr><w:r><w:t>Delete Empty File~</w:t></w:r></w:p>"

$p2 = $row2.Cells.Item(2).Range.Paragraphs.Item(1)
$innerXml2 = "<w:p><w:pPr><w:pStyle w:val='Compact'/></w:pPr><w:r><w:t xml:space='preserve'>Deletes any file that is 0 bytes in size. Read “Empty File </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t>Error” find procedure for more information. (Includes owner column.)</w:t></w:r></w:p>"
Set-ParaXml $p2 $innerXml2

# 3. The following row ("Replace Space w Hyphen (DIR)") gains the
#    lastRenderedPageBreak that used to sit on "Delete Empty File~".
$row3 = $t.Rows.Item(3)
$p3 = $row3.Cells.Item(1).Range.Paragraphs.Item(1)
Set-ParaXml $p3 "<w:p><w:pPr><w:pStyle w:val='Compact'/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Replace Space w Hyphen (DIR)</w:t></w:r></w:p>"

# 4. Add a new, fully empty row at the end of the table.
$newRow = $t.Rows.Add()
$cell1 = $newRow.Cells.Item(1)
$cell2 = $newRow.Cells.Item(2)
$emptyPara = "<w:p><w:pPr><w:pStyle w:val='Compact'/></w:pPr></w:p>"
Set-ParaXml $cell1.Range.Paragraphs.Item(1) $emptyPara
Set-ParaXml $cell2.Range.Paragraphs.Item(1) $emptyPara
